$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the Correspond Handoff/Handback Datetime for the
# b1beacd1...786bfc1b...zh-cn.xlf row (row 7)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D7").Value = "2016-03-08 06:37:52"
$wsZh.Range("G7").Value = "2016-03-08 06:38:33"

# de-de sheet: update the Correspond Handoff/Handback Datetime for the
# b1beacd1...786bfc1b...de-de.xlf row (row 7)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D7").Value = "2016-03-08 06:38:02"
$wsDe.Range("G7").Value = "2016-03-08 06:38:48"
